# Weekly price-log update: prepend a new week's worth of "Apio" (celery)
# price records for "Vega Central Mapocho de Santiago" ahead of the most
# recent existing entries (which currently start at row 442), pushing all
# subsequent historical rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (row 442-443), which
# shifts the former row 442 onward down to rows 444+ and grows the used
# range from R556 to R558 automatically.
$ws.Range("442:443").Insert()

# Row 442 - "Primera" quality record for the new week (2023-10-05).
$ws.Range("A442").Value = 9
$ws.Range("B442").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C442").Value = "Metropolitana"
$ws.Range("D442").Value = 45204
$ws.Range("E442").Value = 13
$ws.Range("F442").Value = 100112017
$ws.Range("G442").Value = "Apio"
$ws.Range("H442").Value = "Americana (o)"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 70
$ws.Range("K442").Value = 6000
$ws.Range("L442").Value = 7000
$ws.Range("M442").Value = 6500
$ws.Range("N442").Value = "$/docena de matas"
$ws.Range("O442").Value = "Región de Coquimbo"
$ws.Range("P442").Value = 1083
$ws.Range("Q442").Value = 6
$ws.Range("R442").Value = "Hortaliza"

# Row 443 - "Segunda" quality record for the new week (2023-10-05).
$ws.Range("A443").Value = 9
$ws.Range("B443").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C443").Value = "Metropolitana"
$ws.Range("D443").Value = 45204
$ws.Range("E443").Value = 13
$ws.Range("F443").Value = 100112017
$ws.Range("G443").Value = "Apio"
$ws.Range("H443").Value = "Americana (o)"
$ws.Range("I443").Value = "Segunda"
$ws.Range("J443").Value = 26
$ws.Range("K443").Value = 5000
$ws.Range("L443").Value = 5000
$ws.Range("M443").Value = 5000
$ws.Range("N443").Value = "$/docena de matas"
$ws.Range("O443").Value = "Región de Coquimbo"
$ws.Range("P443").Value = 833
$ws.Range("Q443").Value = 6
$ws.Range("R443").Value = "Hortaliza"

# Note: Range("442:443").Insert() already carries the row's existing
# formatting (incl. the date number-format on column D) onto the two new
# rows, so no extra style copy is needed here.
